$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.590.32"
$ws.Range("E2").Value = "  -1.37%  "
$ws.Range("D3").Value = "2.031.46"
$ws.Range("E3").Value = "  +1.13%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.16"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -9.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.601"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.90%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "55.08"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.82%  "
$ws.Range("E9").Value = "  -1.94%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.55"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.65%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0752"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.02%  "
$ws.Range("E12").Value = "  -1.72%  "
$ws.Range("D13").Value = "2.330.77"
$ws.Range("E13").Value = "  +1.00%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.25"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.28"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -5.58%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.763"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -4.21%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.09"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -2.05%  "
$ws.Range("D18").Value = "2.029.87"
$ws.Range("E18").Value = "  -3.82%  "
$ws.Range("D19").Value = "36.607.55"
$ws.Range("E19").Value = "  -1.05%  "
$ws.Range("E20").Value = "  -4.36%  "
$ws.Range("E21").Value = "  -3.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.41"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +6.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "221.31"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -5.47%  "
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("E25").Value = "  +1.23%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.40"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -6.55%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.81"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.26%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.134"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +7.23%  "
$ws.Range("E29").Value = "  -3.37%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "18.99"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.73%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.37"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.77%  "
$ws.Range("E32").Value = "  -1.81%  "
$ws.Range("E33").Value = "  -4.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0603"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -6.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.48"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +5.14%  "
$ws.Range("E36").Value = "  -3.36%  "
$ws.Range("E37").Value = "  +0.11%  "
$ws.Range("E38").Value = "  -3.74%  "
$ws.Range("E39").Value = "  -2.68%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.80"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +4.49%  "
$ws.Range("B41").Value = "HuobiToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.91"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -3.98%  "
$ws.Range("B42").Value = "Cronos"
$ws.Range("C42").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0952"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +3.52%  "
$ws.Range("D43").Value = "1.458.66"
$ws.Range("E43").Value = "  +1.35%  "
$ws.Range("E44").Value = "  -2.91%  "
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.11"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -5.33%  "
$ws.Range("B46").Value = "FTXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.11"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +37.71%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.77"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.71%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.60"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.33%  "
$ws.Range("E49").Value = "  -1.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.89"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.58%  "
$ws.Range("E51").Value = "  -1.11%  "
